{"js": "// Apply the \"Built site for gh-pages\" style changes to the AMCR Knihovna 3D\n// template: add the \"Abstract Title\" paragraph style, tighten the spacing\n// above the existing \"Abstract\" style, and add the \"Footnote Block Text\"\n// paragraph style.\n\n// ------------------------------------------------------------------\n// 1) New style: \"Abstract Title\" (w:styleId=\"AbstractTitle\")\n// ------------------------------------------------------------------\ncontext.document.addStyle(\"Abstract Title\", Word.StyleType.paragraph);\nawait context.sync();\n\nconst abstractTitle = context.document.getStyles().getByNameOrNullObject(\"Abstract Title\");\nawait context.sync();\n\nabstractTitle.baseStyle = \"Normal\";\nabstractTitle.nextParagraphStyle = \"Abstract\";\nabstractTitle.quickStyle = true;\n\nabstractTitle.paragraphFormat.keepWithNext = true;\nabstractTitle.paragraphFormat.keepTogether = true;\nabstractTitle.paragraphFormat.alignment = Word.Alignment.centered;\nabstractTitle.paragraphFormat.spaceBefore = 15;\nabstractTitle.paragraphFormat.spaceAfter = 0;\n\nabstractTitle.font.size = 10;\nabstractTitle.font.sizeBidirectional = 10;\nabstractTitle.font.bold = true;\nabstractTitle.font.color = \"#345A8A\";\nawait context.sync();\n\n// ------------------------------------------------------------------\n// 2) Existing style \"Abstract\": tighten the space above the paragraph\n//    (w:before 300 -> 100 twips ; 15pt -> 5pt). SpaceAfter is untouched.\n// ------------------------------------------------------------------\nconst abstractStyle = context.document.getStyles().getByNameOrNullObject(\"Abstract\");\nawait context.sync();\nabstractStyle.paragraphFormat.spaceBefore = 5;\nawait context.sync();\n\n// ------------------------------------------------------------------\n// 3) New style: \"Footnote Block Text\" (w:styleId=\"FootnoteBlockText\")\n// ------------------------------------------------------------------\ncontext.document.addStyle(\"Footnote Block Text\", Word.StyleType.paragraph);\nawait context.sync();\n\nconst footnoteBlockText = context.document.getStyles().getByNameOrNullObject(\"Footnote Block Text\");\nawait context.sync();\n\nfootnoteBlockText.baseStyle = \"Footnote Text\";\nfootnoteBlockText.nextParagraphStyle = \"Footnote Text\";\nfootnoteBlockText.priority = 9;\nfootnoteBlockText.unhideWhenUsed = true;\nfootnoteBlockText.quickStyle = true;\n\nfootnoteBlockText.paragraphFormat.spaceBefore = 5;\nfootnoteBlockText.paragraphFormat.spaceAfter = 5;\nfootnoteBlockText.paragraphFormat.firstLineIndent = 0;\nfootnoteBlockText.paragraphFormat.leftIndent = 24;\nfootnoteBlockText.paragraphFormat.rightIndent = 24;\nawait context.sync();\n", "ps1": "# Apply the \"Built site for gh-pages\" style changes to the AMCR Knihovna 3D\n# template: add the \"Abstract Title\" paragraph style, tighten the spacing\n# above the existing \"Abstract\" style, and add the \"Footnote Block Text\"\n# paragraph style.\n\n$d = $word.ActiveDocument\n\n# ------------------------------------------------------------------\n# 1) New style: \"Abstract Title\" (w:styleId=\"AbstractTitle\")\n# ------------------------------------------------------------------\n$abstractTitle = $d.Styles.Add(\"Abstract Title\", 1)\n$abstractTitle.BaseStyle = $d.Styles(\"Normal\")\n$abstractTitle.NextParagraphStyle = \"Abstract\"\n$abstractTitle.QuickStyle = $true\n\n$abstractTitle.ParagraphFormat.KeepWithNext = $true\n$abstractTitle.ParagraphFormat.KeepTogether = $true\n$abstractTitle.ParagraphFormat.Alignment = 1\n$abstractTitle.ParagraphFormat.SpaceBefore = 15\n$abstractTitle.ParagraphFormat.SpaceAfter = 0\n\n$abstractTitle.Font.Size = 10\n$abstractTitle.Font.SizeBi = 10\n$abstractTitle.Font.Bold = $true\n$abstractTitle.Font.Color = 9067060\n\n# ------------------------------------------------------------------\n# 2) Existing style \"Abstract\": tighten the space above the paragraph\n#    (w:before 300 -> 100 twips ; 15pt -> 5pt). SpaceAfter is untouched.\n# ------------------------------------------------------------------\n$abstract = $d.Styles(\"Abstract\")\n$abstract.ParagraphFormat.SpaceBefore = 5\n\n# ------------------------------------------------------------------\n# 3) New style: \"Footnote Block Text\" (w:styleId=\"FootnoteBlockText\")\n# ------------------------------------------------------------------\n$footnoteBlockText = $d.Styles.Add(\"Footnote Block Text\", 1)\n$footnoteBlockText.BaseStyle = \"Footnote Text\"\n$footnoteBlockText.NextParagraphStyle = \"Footnote Text\"\n$footnoteBlockText.Priority = 9\n$footnoteBlockText.UnhideWhenUsed = $true\n$footnoteBlockText.QuickStyle = $true\n\n$footnoteBlockText.ParagraphFormat.SpaceBefore = 5\n$footnoteBlockText.ParagraphFormat.SpaceAfter = 5\n$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0\n$footnoteBlockText.ParagraphFormat.LeftIndent = 24\n$footnoteBlockText.ParagraphFormat.RightIndent = 24\n"}
